# "Generate Report for Handoff"
# For the 4 files that were "Ready for handoff" (rows 4-7 on the zh-cn and
# de-de sheets), a handoff was generated:
#   - Priority changes from "low" to "ht"
#   - Latest Handoff Datetime is refreshed to the generation time

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-06 00:35:44"

    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-09-06 00:35:49"

    # "Latest HO Xliff Generate Date" on the Overview sheet mirrors the
    # de-de handoff datetime for these rows.
    $overview.Cells.Item($r, 7).Value = "2016-09-06 00:35:49"
}
